# Add season-record columns (Wins / Losses / Ties) to the player stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - copy the style of the existing header cell (AC1) so the
# new headers match the bold/centered/bordered look of the rest of row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Season record values for every player row (2-40).
$lastRow = 40
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 69  # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 93  # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 1   # AF -> Ties
}
